# Applies the "14. feb" update to the "Kjente feil SME 2022" sheet:
# fills in two new known-issue rows (rows 8 and 9) that previously were
# blank placeholder rows in the table, and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kjente feil SME 2022")

# Row 8 - issue #3
$ws.Range("A8").Value = 3
$ws.Range("B8").Value = (Get-Date -Year 2023 -Month 2 -Day 14 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("C8").Value = "Etter å ha oppettet en instans av skattemelding i altinn kan man ikke oppdatere metadata på instansen (med for eksempel at signeres av revisor). Det må settes i den initielle opprettelsen. Det jobbes med feilfiks."
$ws.Range("D8").Value = "Altinn"
$ws.Range("E8").Value = "AS"
$ws.Range("F8").Value = "Produksjon"
$ws.Range("G8").Value = "Åpen"

# Row 9 - issue #4
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = (Get-Date -Year 2023 -Month 2 -Day 14 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("C9").Value = "i informasjonsmodellen for FOU står det at 'Vedtak', 'AvgrensingIVedtak' at denne kan være 0 til mange. Dette er feil - det vil kun komme en forekomst her. Vi endrer ikke informasjonsmodellen for dette nå i år. Sender dere inn flere forekomster, som er lovlige iht XSD vil kun første forekomst faktisk bli sendt inn. "
$ws.Range("D9").Value = "Altinn"
$ws.Range("E9").Value = "AS"
$ws.Range("F9").Value = "Produksjon"
$ws.Range("G9").Value = "Åpen"

# Row 9 grew to fit its long wrapped description - match the autofitted height.
$ws.Rows.Item(9).RowHeight = 48

# Move the active selection to A9, matching the saved view state.
$ws.Range("A9").Select()
